$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts following the rescaling of the MALI clustering (incorporating time)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
